# Commit: Add indst/ItUBB to acronym key (commit #fe63694; 8/12/24)
#
# In the "Key to Variables" sheet, a new row is inserted for the "indst"
# (Industrial) top-level folder, introducing the acronym "ItUBB"
# ("Industries that Use Byproduct Biomass"). The new row is inserted
# right before the existing "indst" / "MHV" row (sheet row 188 after the
# insert; row 187 before it), pushing that row and everything below it
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a new blank row immediately above the current row 187
# ("indst" / "MHV" / "Methane Heating Value"), shifting it (and all
# subsequent rows) down to row 188.
$ws.Rows.Item(187).Insert()

# The row-insert copies the row-above's used-column formatting, which
# drags a spurious formatted cell into column G (row 186 has data in G).
# The new row has no "Update Only Needed If" note, so drop it entirely.
$ws.Range("G187").Clear()

# Populate the new row's data:
#   A = Top Level Folder -> indst (Industrial)
#   B = Acronym           -> ItUBB
#   C = Meaning            -> Industries that Use Byproduct Biomass
#   F = Importance to Update for New Country -> low
$ws.Range("A187").Value = "indst"
$ws.Range("B187").Value = "ItUBB"
$ws.Range("C187").Value = "Industries that Use Byproduct Biomass"
$ws.Range("F187").Value = "low"

# Match the "low"-importance cell shading used by sibling rows (e.g. F185)
# rather than the style inherited from the row above.
$ws.Range("F185").Copy()
$ws.Range("F187").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to reflect where the editor ended up after
# making the edit.
$ws.Range("D190").Select()
